# Update the cryptos list with latest scraped prices / 1h volume changes.
# Generated from the GitHub Actions commit: "Updated cryptos list on Wed Jan 17 03:42:33 UTC 2024 with GitHub Actions"
#
# NOTE: the "Price" column (D) stores values as plain text (e.g. "42.953.72"),
# not numbers. Values that look like a single decimal number (e.g. "315.42")
# would otherwise be auto-converted to a real number by Excel, which both
# changes their type and can silently drop significant trailing zeros
# (e.g. "36.30" -> 36.3). Prefixing those with a leading apostrophe forces
# Excel to keep them as literal text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.911.34"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.575.19"
$ws.Range("E3").Value = "  +1.95%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'315.42"
$ws.Range("E5").Value = "  +0.32%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'99.63"
$ws.Range("E6").Value = "  +3.98%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.18%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "  +0.49%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'36.30"
$ws.Range("E10").Value = "  +0.38%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.44%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "'7.51"
$ws.Range("E12").Value = "  -0.55%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.972.93"
$ws.Range("E13").Value = "  +2.14%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.56%  "

# Rows 15 & 16 - Chainlink and WrappedEther swap ranking positions
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.641.52"
$ws.Range("E15").Value = "  +6.56%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'15.73"
$ws.Range("E16").Value = "  +2.90%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -1.07%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.982.33"
$ws.Range("E18").Value = "  +0.25%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +1.35%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = "  -2.00%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0968"
$ws.Range("E21").Value = "  +0.47%  "

# Row 22 - Litecoin
$ws.Range("E22").Value = "  -0.30%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'249.95"
$ws.Range("E23").Value = "  -1.62%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "'2.97"
$ws.Range("E24").Value = "  +0.55%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -0.07%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "'27.12"
$ws.Range("E26").Value = "  +1.74%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.00%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  -1.20%  "

# Row 29 - InjectiveProtocol
$ws.Range("D29").Value = "'40.66"
$ws.Range("E29").Value = "  -0.55%  "

# Row 30 - Cosmos
$ws.Range("D30").Value = "'10.31"
$ws.Range("E30").Value = "  -0.73%  "

# Row 31 - Monero
$ws.Range("D31").Value = "'158.05"
$ws.Range("E31").Value = "  +0.32%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -1.75%  "

# Row 33 - LidoDAOToken
$ws.Range("D33").Value = "'3.44"
$ws.Range("E33").Value = "  +3.77%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  -1.54%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "'0.0803"
$ws.Range("E35").Value = "  +3.05%  "

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  -0.33%  "

# Row 37 - Celestia
$ws.Range("D37").Value = "'18.86"
$ws.Range("E37").Value = "  -2.80%  "

# Row 38 - ApeXProtocol
$ws.Range("E38").Value = "  +9.82%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +1.24%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +0.29%  "

# Row 41 - EnergySwap
$ws.Range("D41").Value = "'23.59"
$ws.Range("E41").Value = "  -0.01%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "'4.09"
$ws.Range("E42").Value = "  +7.81%  "

# Row 43 - VeChain
$ws.Range("D43").Value = "'0.0304"
$ws.Range("E43").Value = "  -0.40%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  -0.06%  "

# Row 45 - NEARProtocol
$ws.Range("D45").Value = "'3.25"
$ws.Range("E45").Value = "  -2.40%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.005.71"
$ws.Range("E46").Value = "  -2.39%  "

# Row 47 - FraxShare
$ws.Range("D47").Value = "'8.94"
$ws.Range("E47").Value = "  +0.21%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "2.822.62"
$ws.Range("E48").Value = "  +2.08%  "

# Row 49 - Algorand
$ws.Range("D49").Value = "'0.197"
$ws.Range("E49").Value = "  +2.71%  "

# Row 50 - ordi
$ws.Range("D50").Value = "'75.11"
$ws.Range("E50").Value = "  -0.49%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "'81.76"
$ws.Range("E51").Value = "  -4.40%  "
